{"js": "// Convert the \"{ m:self.name }\" Word FIELD (fldChar begin/instrText/fldChar end)\n// into plain literal text runs \"{\", \"m\", \":\", \"self\", \".name}\" - i.e. turn the\n// field code into literal M2Doc template text, keeping the orange highlight\n// that was on the \"self\" token.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that actually contains the field (robust to position).\nfor (const p of paragraphs.items) {\n  p.fields.load(\"items\");\n}\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.fields.items.length > 0) {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  // Recover the original paragraph-level attributes (w:rsidR/.../w:rsidP) so\n  // the rewritten paragraph keeps looking like the same paragraph, not a\n  // brand-new one.\n  const originalOoxml = target.getOoxml();\n  await context.sync();\n\n  let pAttrs = \"\";\n  const m = originalOoxml.value.match(/<w:p\\b([^>]*)>/);\n  if (m) {\n    // Drop the synthetic w14:paraId/w14:textId that getOoxml() stamps on -\n    // they are not part of the original part.\n    pAttrs = m[1]\n      .replace(/\\s*w14:paraId=\"[^\"]*\"/, \"\")\n      .replace(/\\s*w14:textId=\"[^\"]*\"/, \"\");\n  }\n\n  // Build the replacement paragraph as literal text runs. The \"self\" run\n  // keeps the original orange accent color (incl. theme info) that used to\n  // live on its <w:rPr>; the other runs become plain text runs. The field\n  // wrapper characters \"{\" and \"}\" are added around the (now literal) field\n  // code, and the inner whitespace that used to pad the field code\n  // (\"  m:self.name  \") is trimmed away.\n  const replacementOoxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    \"<w:p\" + pAttrs + \">\" +\n    \"<w:r><w:t>{</w:t></w:r>\" +\n    \"<w:r><w:t>m</w:t></w:r>\" +\n    \"<w:r><w:t>:</w:t></w:r>\" +\n    '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>self</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">.name}</w:t></w:r>' +\n    \"</w:p>\" +\n    \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\";\n\n  target.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Convert the \"{ m:self.name }\" Word FIELD (fldChar begin/instrText/fldChar end)\n# into plain literal text runs \"{\", \"m\", \":\", \"self\", \".name}\" - i.e. turn the\n# field code into literal M2Doc template text, keeping the orange highlight\n# that was on the \"self\" token.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that actually contains the field (robust to position).\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Fields.Count -gt 0) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $rng = $target.Range\n\n    # Recover the original paragraph-level attributes (w:rsidR/.../w:rsidP) so\n    # the rewritten paragraph keeps looking like the same paragraph, not a\n    # brand-new one.\n    $origXml = $rng.XML()\n    $pAttrs = \"\"\n    if ($origXml -match '<w:p\\b([^>]*)>') {\n        $pAttrs = $matches[1]\n        # Drop the synthetic w14:paraId/w14:textId that Range.XML() stamps on -\n        # they are not part of the original part.\n        $pAttrs = $pAttrs -replace ' w14:paraId=\"[^\"]*\"', ''\n        $pAttrs = $pAttrs -replace ' w14:textId=\"[^\"]*\"', ''\n    }\n\n    # Build the replacement paragraph as literal text runs. The \"self\" run\n    # keeps the original orange accent color (incl. theme info) that used to\n    # live on its <w:rPr>; the other runs become plain text runs. The field\n    # wrapper characters \"{\" and \"}\" are added around the (now literal) field\n    # code, and the inner whitespace that used to pad the field code\n    # (\"  m:self.name  \") is trimmed away.\n    $replacementXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n        '<w:p' + $pAttrs + '>' +\n        '<w:r><w:t>{</w:t></w:r>' +\n        '<w:r><w:t>m</w:t></w:r>' +\n        '<w:r><w:t>:</w:t></w:r>' +\n        '<w:r><w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr><w:t>self</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\">.name}</w:t></w:r>' +\n        '</w:p>' +\n        '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n\n    $rng.InsertXML($replacementXml)\n}\n"}
